# StructureDefinition-enrollment-pcp.xlsx — FHIR version bump + publisher/jurisdiction update
# (Alvearie alvearie-fhir-ig gh-pages deploy: 5.0.0 -> 6.0.0, Contact -> Publisher/Jurisdiction,
#  and the root Extension's Short/Definition filled in on the Elements sheet.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet (sheet1): property/value pairs in columns A/B
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: regenerated on this publish
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty; now populated
$meta.Range("B9").Value = "Alvearie Team"

# The old sheet had a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11). Row 10 becomes the new "Jurisdiction" row, and the
# duplicate row 11 is removed entirely (shifting everything below up by one).
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# "Elements" sheet (sheet2): the root Extension row (row 2) gets a real
# Short/Definition instead of the generic "Extension" / "An Extension"
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Enrollment PCP"
$elements.Range("L2").Value = "Provider identifier of the primary care physician for the plan member"
